# Apply cached-value corrections to the profit-calculation columns (H:N)
# across several leve rows on multiple sheets, per the upstream data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 386.44186
$ws.Range("J17").Value = 390.88095
$ws.Range("L17").Value = 1172.64285
$ws.Range("N17").Value = -1508.64285

$ws.Range("H28").Value = 1692.7858
$ws.Range("I28").Value = 1719.52
$ws.Range("J28").Value = 1470
$ws.Range("K28").Value = 1719.52
$ws.Range("L28").Value = 1470
$ws.Range("M28").Value = -1234.52
$ws.Range("N28").Value = -2440

$ws.Range("H40").Value = 3790.25
$ws.Range("I40").Value = 3450
$ws.Range("J40").Value = 4130.5
$ws.Range("K40").Value = 3450
$ws.Range("L40").Value = 4130.5
$ws.Range("M40").Value = -3275
$ws.Range("N40").Value = -4480.5

$ws.Range("H62").Value = 5966.6665
$ws.Range("I62").Value = 5966.6665
$ws.Range("K62").Value = 5966.6665
$ws.Range("M62").Value = -5342.6665

$ws.Range("H65").Value = 5966.6665
$ws.Range("I65").Value = 5966.6665
$ws.Range("K65").Value = 29833.3325
$ws.Range("M65").Value = -26713.3325

$ws.Range("H113").Value = 4991.55
$ws.Range("I113").Value = 5463.7334
$ws.Range("J113").Value = 3575
$ws.Range("K113").Value = 5463.7334
$ws.Range("L113").Value = 3575
$ws.Range("M113").Value = -2209.7334
$ws.Range("N113").Value = -10083

$ws.Range("H127").Value = 1921681.2
$ws.Range("I127").Value = 1584.2307
$ws.Range("K127").Value = 4752.6921
$ws.Range("M127").Value = 207.3078999999998

$ws.Range("H132").Value = 4995.7314
$ws.Range("I132").Value = 4067.2856
$ws.Range("J132").Value = 7523.1665
$ws.Range("K132").Value = 12201.8568
$ws.Range("L132").Value = 22569.4995
$ws.Range("M132").Value = -9671.856800000001
$ws.Range("N132").Value = -27629.4995

$ws.Range("H138").Value = 21778696
$ws.Range("J138").Value = 2163477
$ws.Range("L138").Value = 6490431
$ws.Range("N138").Value = -6500711

$ws.Range("H140").Value = 92842.30499999999
$ws.Range("J140").Value = 92842.30499999999
$ws.Range("L140").Value = 92842.30499999999
$ws.Range("N140").Value = -103202.305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4361.3594
$ws.Range("I32").Value = 3486.3
$ws.Range("K32").Value = 3486.3
$ws.Range("M32").Value = -3199.3

$ws.Range("H45").Value = 6024.727
$ws.Range("I45").Value = 8053.45
$ws.Range("J45").Value = 2903.6155
$ws.Range("K45").Value = 8053.45
$ws.Range("L45").Value = 2903.6155
$ws.Range("M45").Value = -7676.45
$ws.Range("N45").Value = -3657.6155

$ws.Range("H132").Value = 2750.3914
$ws.Range("I132").Value = 2108.1
$ws.Range("J132").Value = 7032.3335
$ws.Range("K132").Value = 6324.299999999999
$ws.Range("L132").Value = 21097.0005
$ws.Range("M132").Value = -3794.299999999999
$ws.Range("N132").Value = -26157.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3630.6765
$ws.Range("I134").Value = 3610.6072
$ws.Range("J134").Value = 3724.3333
$ws.Range("K134").Value = 10831.8216
$ws.Range("L134").Value = 11172.9999
$ws.Range("M134").Value = -8296.821599999999
$ws.Range("N134").Value = -16242.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3879.8147
$ws.Range("I31").Value = 2447.5833
$ws.Range("J31").Value = 5025.6
$ws.Range("K31").Value = 2447.5833
$ws.Range("L31").Value = 5025.6
$ws.Range("M31").Value = -2152.5833
$ws.Range("N31").Value = -5615.6

$ws.Range("H34").Value = 3879.8147
$ws.Range("I34").Value = 2447.5833
$ws.Range("J34").Value = 5025.6
$ws.Range("K34").Value = 2447.5833
$ws.Range("L34").Value = 5025.6
$ws.Range("M34").Value = -2245.5833
$ws.Range("N34").Value = -5429.6

$ws.Range("H51").Value = 39065
$ws.Range("J51").Value = 69498.75
$ws.Range("L51").Value = 69498.75
$ws.Range("N51").Value = -70970.75

$ws.Range("H61").Value = 39065
$ws.Range("J61").Value = 69498.75
$ws.Range("L61").Value = 69498.75
$ws.Range("N61").Value = -70194.75

$ws.Range("H107").Value = 1341.7812
$ws.Range("I107").Value = 558.9524
$ws.Range("J107").Value = 2836.2727
$ws.Range("K107").Value = 558.9524
$ws.Range("L107").Value = 2836.2727
$ws.Range("M107").Value = 1361.0476
$ws.Range("N107").Value = -6676.2727

$ws.Range("H122").Value = 1950.2258
$ws.Range("I122").Value = 1612.84
$ws.Range("J122").Value = 3356
$ws.Range("K122").Value = 4838.52
$ws.Range("L122").Value = 10068
$ws.Range("M122").Value = -2388.52
$ws.Range("N122").Value = -14968

$ws.Range("H134").Value = 4152.017
$ws.Range("I134").Value = 3211.3408
$ws.Range("J134").Value = 7108.4287
$ws.Range("K134").Value = 9634.0224
$ws.Range("L134").Value = 21325.2861
$ws.Range("M134").Value = -7099.0224
$ws.Range("N134").Value = -26395.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 25046000
$ws.Range("I32").Value = 7500
$ws.Range("J32").Value = 33392166
$ws.Range("K32").Value = 22500
$ws.Range("L32").Value = 100176498
$ws.Range("M32").Value = -22217
$ws.Range("N32").Value = -100177064

$ws.Range("H68").Value = 19158.834
$ws.Range("J68").Value = 22590.6
$ws.Range("L68").Value = 67771.79999999999
$ws.Range("N68").Value = -69393.79999999999

$ws.Range("H71").Value = 19158.834
$ws.Range("J71").Value = 22590.6
$ws.Range("L71").Value = 203315.4
$ws.Range("N71").Value = -211427.4

$ws.Range("H114").Value = 5557110.5
$ws.Range("J114").Value = 8335148
$ws.Range("L114").Value = 25005444
$ws.Range("N114").Value = -25011952

$ws.Range("H126").Value = 16030
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H131").Value = 7719766
$ws.Range("I131").Value = 13370077
$ws.Range("J131").Value = 5298203.5
$ws.Range("K131").Value = 40110231
$ws.Range("L131").Value = 15894610.5
$ws.Range("M131").Value = -40105191
$ws.Range("N131").Value = -15904690.5

$ws.Range("H141").Value = 22117
$ws.Range("I141").Value = 22117
$ws.Range("K141").Value = 66351
$ws.Range("M141").Value = -61171

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3020.2856
$ws.Range("I16").Value = 2857
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 2857
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -2687
$ws.Range("N16").Value = -4340

$ws.Range("H100").Value = 1019.5417
$ws.Range("I100").Value = 885.13635
$ws.Range("K100").Value = 885.13635
$ws.Range("M100").Value = -344.13635

$ws.Range("H122").Value = 5364.154
$ws.Range("I122").Value = 3629.25
$ws.Range("J122").Value = 8140
$ws.Range("K122").Value = 10887.75
$ws.Range("L122").Value = 24420
$ws.Range("M122").Value = -8437.75
$ws.Range("N122").Value = -29320

$ws.Range("H136").Value = 4620827
$ws.Range("I136").Value = 8183422.5
$ws.Range("J136").Value = 10409.059
$ws.Range("K136").Value = 24550267.5
$ws.Range("L136").Value = 31227.177
$ws.Range("M136").Value = -24547717.5
$ws.Range("N136").Value = -36327.177

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9700
$ws.Range("I62").Value = 9000
$ws.Range("J62").Value = 10400
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 10400
$ws.Range("M62").Value = -8376
$ws.Range("N62").Value = -11648

$ws.Range("H65").Value = 9700
$ws.Range("I65").Value = 9000
$ws.Range("J65").Value = 10400
$ws.Range("K65").Value = 45000
$ws.Range("L65").Value = 52000
$ws.Range("M65").Value = -41880
$ws.Range("N65").Value = -58240

$ws.Range("H107").Value = 1323.9231
$ws.Range("I107").Value = 1381.3334
$ws.Range("J107").Value = 1194.75
$ws.Range("K107").Value = 4144.0002
$ws.Range("L107").Value = 3584.25
$ws.Range("M107").Value = -2224.0002
$ws.Range("N107").Value = -7424.25

$ws.Range("H113").Value = 4379
$ws.Range("I113").Value = 4359
$ws.Range("J113").Value = 4419
$ws.Range("K113").Value = 13077
$ws.Range("L113").Value = 13257
$ws.Range("M113").Value = -10907
$ws.Range("N113").Value = -17597

$ws.Range("H122").Value = 6530.676
$ws.Range("I122").Value = 3864.1875
$ws.Range("J122").Value = 23596.2
$ws.Range("K122").Value = 11592.5625
$ws.Range("L122").Value = 70788.60000000001
$ws.Range("M122").Value = -9142.5625
$ws.Range("N122").Value = -75688.60000000001

$ws.Range("H126").Value = 7096.561
$ws.Range("I126").Value = 6753.353
$ws.Range("K126").Value = 20260.059
$ws.Range("M126").Value = -17790.059

$ws.Range("H132").Value = 1652.5385
$ws.Range("I132").Value = 1652.5385
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4957.6155
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2427.6155
$ws.Range("N132").ClearContents()
